# Commit: "adding data structure for evidence"
#
# Insert a new "Evidence" worksheet immediately before the existing
# "References" worksheet (so the tab order becomes:
#   ... Stop conditions, Evidence, References)
# and populate its header row with the columns used to capture
# experimental evidence supporting model values (Taxon, Genetic variant,
# Temperature (C), pH, Growth media, ...).

$wb = $excel.ActiveWorkbook

# Locate the existing "References" sheet and add the new sheet right
# before it, so it lands in the correct tab position.
$refSheet = $wb.Worksheets.Item("References")
$evidence = $wb.Worksheets.Add($refSheet)
$evidence.Name = "Evidence"

# Header row for the new Evidence sheet.
$evidence.Range("A1").Value = "Id"
$evidence.Range("B1").Value = "Name"
$evidence.Range("C1").Value = "Value"
$evidence.Range("D1").Value = "Units"
$evidence.Range("E1").Value = "Type"
$evidence.Range("F1").Value = "Taxon"
$evidence.Range("G1").Value = "Genetic variant"
$evidence.Range("H1").Value = "Temperature (C)"
$evidence.Range("I1").Value = "pH"
$evidence.Range("J1").Value = "Growth media"
$evidence.Range("K1").Value = "Database references"
$evidence.Range("L1").Value = "Evidence"
$evidence.Range("M1").Value = "Comments"
$evidence.Range("N1").Value = "References"

# Match the freeze-pane / selection convention used by the sibling data
# sheets in this workbook (header row frozen, first column frozen).
$evidence.Activate()
$evidence.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
